$wb = $excel.ActiveWorkbook

# --- "feature requirement" sheet: selection moves from C10 to C13 (tab no longer selected) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C13").Select()

# --- "task" sheet: renamed to "route analyst", becomes the active/selected tab ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "route analyst"

# existing header cell
$ws3.Range("A1").Value = "route"

# new route listing content
$ws3.Range("A3").Value = "index"
$ws3.Range("A4").Value = "/api/"
$ws3.Range("B5").Value = "/api-student/"
$ws3.Range("B6").Value = "/api-teacher/"
$ws3.Range("B7").Value = "/api/admin/"
$ws3.Range("A10").Value = "/class-teacher/"
$ws3.Range("A13").Value = "/class-student/"
$ws3.Range("A19").Value = "/test-teacher/"
$ws3.Range("A20").Value = "/test-student/"
$ws3.Range("A23").Value = "/me"
$ws3.Range("A25").Value = "/site"
$ws3.Range("A27").Value = "/chat"

# column widths: A=14.125, B:I=15.5 (closest representable values in this engine's width grid)
$ws3.Columns.Item(1).ColumnWidth = 13.333333333333334
$ws3.Range("B1:I1").EntireColumn.ColumnWidth = 14.666666666666666

# activate the sheet + set the new selection
$ws3.Activate()
$ws3.Range("E9").Select()
